$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Pipette" - add new issue row (row 7)
# ---------------------------------------------------------------------
$wsPipette = $wb.Worksheets.Item("Pipette")
$wsPipette.Activate()

$wsPipette.Range("C7").Value = 43158
$wsPipette.Range("D7").Value = "V1.0"
$wsPipette.Range("E7").Value = "Review"
$wsPipette.Range("F7").Value = "Open"
$wsPipette.Range("G7").Value = "외부 GND 연결용 Jack 추가 - ear jack type"

$wsPipette.Range("G14").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Transformer" - fill in rows 5 and 6 with new issues
# ---------------------------------------------------------------------
$wsTransformer = $wb.Worksheets.Item("Transformer")
$wsTransformer.Activate()

# Row 5
$wsTransformer.Range("C5").Value = 43153
$wsTransformer.Range("D5").Value = "V1.0"
$wsTransformer.Range("E5").Value = "SCH"
$wsTransformer.Range("F5").Value = "Open"
$wsTransformer.Range("G5").Value = "CTX2106xx series의 2차측 구조가 datasheet와 틀림`n- 2차측이 2개로 나눠져 있는것으로 보임"
$wsTransformer.Range("G5").WrapText = $true

# Row 6
$wsTransformer.Range("C6").Value = 43158
$wsTransformer.Range("D6").Value = "V1.0 -B,C"
$wsTransformer.Range("E6").Value = "SCH"
$wsTransformer.Range("F6").Value = "Open"
$wsTransformer.Range("G6").Value = "Plasma 발생 안함`nGas 주입 상태에서도 발생 안함"
$wsTransformer.Range("G6").WrapText = $true
$wsTransformer.Range("H6").Value = "전압을 올려보며 시험 필요"

$wsTransformer.Rows("5:6").RowHeight = 33

$wsTransformer.Range("G16").Select() | Out-Null
